$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column B
$ws.Columns.Item(2).ColumnWidth = 94.85546875

$data = @(
  @("Após a Segmentação da imagem vem...", "a rotulação, que permite realizar a contagem de regiões e objetos. A cada região estabelecida, a cada bloco ou a cada componente conexo, um valor é atribuído. "),
  @("Para objetos em movimento, a importancia de extrair características...", "Que sejam independentes e invariantes em relação à posição. Essas características devem ser robustas a translações, rotações e deformações, já que as imagens podem ter sido capturadas em momentos diferentes."),
  @("Categorização das características diversas utilizadas no processo de reconhecimento ou identificação de objetos em imagens.", "Possui dois descritores: Forma e Aspecto."),
  @("Quando falamos de Forma", "Consideramos apenas o contorno dos objetos e as regiões internas"),
  @("A definição da caracteristica Dimensional, em relação a Regiões.", "Seria a área, perímetro, raio, compacidade etc."),
  @("A definição da caracteristica de Inercia, em relação a Regiões.", "Seria o centro geométrico, orientação, retângulo envolvente etc."),
  @("A definição da caracteristica Topológica, em relação a Regiões.", "Seria componentes conectados, número de vértices, número de Euler etc."),
  @("Definição de rugosidade, em relação a características de aspecto.", "Expressa o acabamento de um objeto e ajuda a detectar regularidades em superfícies"),
  @("Definição de textura, em relação a características de aspecto.", "Pode ser utilizada para caracterizar regularidades como sulcos ou estrias, além disso fornece padronagens de determinados tipos de material, como madeiras, rochas, tecidos etc."),
  @("Definição da Análise de Componentes Principais (ACP).", "Transforma variáveis discretas em coeficientes descorrelacionados, o que permite identificar as relações entre as caraterísticas extraídas de dados."),
  @("Objetivo da Análise de Componentes Principais (ACP).", "Este método tem por objetivo analisar o conjunto de características usado, visando reduzi-lo, a fim de eliminar a redundância."),
  @("Definição de processo autônomo , em relação a classificação de uma imagem.", "É um processo de classificação automático, porém em um primeiro momento é necessário regular manualmente para classificação. necessitando de um humano para fornecer pré-definições possíveis."),
  @("Diferença entre a classificação supervisionada e não supervisionada.", "A classificação não supervisionada agrupa dados sem exemplos prévios, enquanto a classificação supervisionada usa referências conhecidas para classificar novos dados. "),
  @("Objetivo da Definição de classes, para um sistema de classificação em sua primeira etapa.", "Realiza a definição de que tipo de característica o problema trata, qual é a melhor maneira, a priori, de representar esta característica, ou seja, o problema é explicitado. "),
  @("Objetivo da Análise das amostras de dados, para um sistema de classificação em sua segunda etapa.", "Busca identificar variáveis que permitem a discriminação do exemplo dado."),
  @("Objetivo da Definição do espaço de características, para um sistema de classificação em sua terceira etapa.", "Seleciona variáveis apropriadas e forma vetores de características para os dados."),
  @("Objetivo da Construção do classificador, para um sistema de classificação em sua quarta etapa.", "Define superfícies de decisão no espaço de características e constrói o classificador que usa dados de treinamento."),
  @("Objetivo da Testagem do classificador, para um sistema de classificação em sua quinta etapa.", "Busca encontrar novas amostras e classificá-las."),
  @("Objetivo da Validação, para um sistema de classificação em sua sexta etapa.", "Pretende verificar se o desempenho da classificação foi suficientemente bom. Caso não seja satisfatório o resultado, retorna a etapa de definição das classes."),
  @("Definição dos Métodos de reconhecimento por decisão teórica, segundo os métodos de reconhecimento de imagem. (Gonzalez e Woods (2000)).", "Como primeiro método, ele faz referência à representação dos padrões na forma de vetor, seguido da busca por abordagens para agrupamento e para atribuição dos vetores de padrões a classes de padrões."),
  @("Definição dos Métodos estruturais para o reconhecimento, segundo os métodos de reconhecimento de imagem. (Gonzalez e Woods (2000)).", "Como segundo método, fazem reconhecimentos estruturais, buscando apresentar os padrões de uma representação simbólica, baseiam-se em casamento de cadeias ou em modelos que tratam os padrões simbólicos como estruturas de uma linguagem artificial."),
  @("Definição dos Métodos para a interpretação de imagens, segundo os métodos de reconhecimento de imagem. (Gonzalez e Woods (2000)).", "Como terceiro método, interpreta imagens, de modo a realizar uma atribuição de significado ao conjunto de elementos de imagens reconhecidas."),
  @("Existem dois métodos em particular que calculam e determinam vetores que representarão o padrão ou as classes de padrões de decisão.", "Esse trabalho são feitos por classificadores,  que pode considerar a praticidade de implementação, como os classificadores pautados em distâncias; ou considerar algum elemento estatístico, que toma por base a utilização de amostras de padrões para estimar determinados parâmetros estatísticos de cada classe de padrões."),
  @("Definição dos Métodos estruturais.", "Procuram realizar reconhecimento de padrão realizando relações estruturais inerentes à forma do padrão."),
  @("Definição de Casamento de números de formas, em relação a Métodos Estruturais", "Refere-se à correspondência entre padrões geométricos (formas) extraídos de uma imagem e um conjunto de primitivas de padrões previamente definido."),
  @("Definição de Casamento de cadeias, em relação a Métodos Estruturais.", "Envolve encontrar correspondências entre padrões (ou cadeias de caracteres) em uma imagem e um modelo previamente definido."),
  @("Definição de Métodos Sintáticos, em relação a Métodos Estruturais.", "Visam ao reconhecimento sintático de padrões que estão no conjunto de primitivas de padrões que, de modo geral, são um conjunto de regras que governa suas interconexões, e um reconhecedor cuja estrutura é determinada pelo conjunto de regras da gramática."),
  @("Definição de Limitação da generalidade do problema, uma das técnicas que podem ser aplicadas na interpretação de uma imagem.", "Se refere ao fato de que, embora tenhamos ferramentas automatizadas para analisar características visuais em imagens, a interpretação crítica ainda depende da percepção humana."),
  @("Definição de Grafos, uma das técnicas que podem ser aplicadas na interpretação de uma imagem.", "Podem ser usados para encontrar padrões, relacionar objetos ou descobrir caminhos relevantes. Por exemplo, se quisermos identificar objetos em uma foto, podemos criar um grafo onde cada objeto é um vértice e as conexões representam relações entre eles."),
  @("Definição de Incorporação de conhecimento humano, uma das técnicas que podem ser aplicadas na interpretação de uma imagem.", "Refere-se a como os seres humanos contribuem com seu entendimento e expertise para melhorar a análise e compreensão de imagens. "),
  @("Definição de Sistemas de produções, uma das técnicas que podem ser aplicadas na interpretação de uma imagem.", "Seria uma linha de montagem para entender e analisar visualmente o que vemos. Eles nos ajudam a “fabricar” significados a partir das imagens que encontramos. Comumente chamada de sistemas baseados em regras."),
  @("Definição de Padrão. (Gonzalez e Woods (2000))", "É uma descrição quantitativa ou estrutural de um objeto ou alguma outra entidade de interesse em uma imagem. "),
  @("Definição de Classe de Padrões", "É um conjunto, uma família de padrões que possui algumas propriedades em comum."),
  @("Um exemplo de relações estruturais é a estrutura em escada.", "É uma estrutura que visa explorar as relações entre suas regiões."),
  @("Em relação a descrição de objetos em imagens, as caracteristicas de Representações por Cadeias são:", "Usam conexões simples entre primitivas (partes) dos objetos. Essas conexões geralmente estão relacionadas à borda externa da forma."),
  @("Em relação a descrição de objetos em imagens, as caracteristicas de Descrição por Árvores são:", "É mais robusta e detalhada e funciona bem para objetos organizados hierarquicamente (como partes de um todo)."),
  @("Definição de Redes Neurais", "São caracterizadas pelo desenvolvimento adaptativo dos coeficientes das funções de decisão através de apresentações sucessivas de conjuntos de padrões de treinamento.")
)

$startRow = 43
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Copy formatting (style) from row 42 to the new rows 43-79, matching column A/B style index 3
$srcRange = $ws.Range("A42:B42")
$dstRange = $ws.Range("A43:B79")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats

# Autofit all rows so wrapped text heights recompute
$ws.Rows.Item("1:79").AutoFit()

# Update view: scroll to bottom, select A80
$ws.Application.ActiveWindow.ScrollRow = 70
$ws.Range("A80").Select()
